$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsBirds = $wb.Worksheets.Item("Birds")
$wsCages = $wb.Worksheets.Item("Cages")

# --- Order matters here: new shared-string values must be written for the
# --- first time in this exact sequence: "arvili75", "a42", "a43".

# 1) Users sheet - new row 34 for the new user "arvili75"
$wsUsers.Range("A34").Value = "arvili75"
$wsUsers.Range("B34").Value = "123456a!"
$wsUsers.Range("C34").Value = 8034

# 2) Birds sheet - new rows 54-56
# Prime the date formatting (column G) by copying the style from the row
# above (G53), which already uses the short-date cell style, then assign
# the actual date serials afterwards so we don't create a brand-new
# number format.
$wsBirds.Range("G53").Copy()
$wsBirds.Range("G54").PasteSpecial(-4122)
$wsBirds.Range("G55").PasteSpecial(-4122)
$wsBirds.Range("G56").PasteSpecial(-4122)

# Row 55 references cage "a42" first, establishing it before "a43" is used
# in row 54, so write H55 before H54 to match the original authoring order.
$wsBirds.Range("H55").Value = "a42"
$wsBirds.Range("H54").Value = "a43"
$wsBirds.Range("H56").Value = "a42"

$wsBirds.Range("A54").Value = 53
$wsBirds.Range("B54").Value = "American Gouldian"
$wsBirds.Range("C54").Value = "Notrh America"
$wsBirds.Range("D54").Value = "Male"
$wsBirds.Range("G54").Value = 45077
$wsBirds.Range("I54").Value = 8034
$wsBirds.Range("J54").Value = "Red"
$wsBirds.Range("K54").Value = "Purple"
$wsBirds.Range("L54").Value = "Green"

$wsBirds.Range("A55").Value = 54
$wsBirds.Range("B55").Value = "American Gouldian"
$wsBirds.Range("C55").Value = "Notrh America"
$wsBirds.Range("D55").Value = "Female"
$wsBirds.Range("G55").Value = 45077
$wsBirds.Range("I55").Value = 8034
$wsBirds.Range("J55").Value = "Red"
$wsBirds.Range("K55").Value = "Purple"
$wsBirds.Range("L55").Value = "Green"

$wsBirds.Range("A56").Value = 55
$wsBirds.Range("B56").Value = "American Gouldian"
$wsBirds.Range("C56").Value = "Notrh America"
$wsBirds.Range("D56").Value = "Male"
$wsBirds.Range("E56").Value = 54
$wsBirds.Range("F56").Value = 53
$wsBirds.Range("G56").Value = 45077
$wsBirds.Range("I56").Value = 8034
$wsBirds.Range("J56").Value = "Red"
$wsBirds.Range("K56").Value = "Purple"
$wsBirds.Range("L56").Value = "Green"

# 3) Cages sheet - new rows 42-43 (reuse the "a42"/"a43" strings already
# registered above)
$wsCages.Range("A42").Value = "a42"
$wsCages.Range("B42").Value = 2
$wsCages.Range("C42").Value = 2
$wsCages.Range("D42").Value = 144
$wsCages.Range("E42").Value = "wood"
$wsCages.Range("F42").Value = 8034

$wsCages.Range("A43").Value = "a43"
$wsCages.Range("B43").Value = 2
$wsCages.Range("C43").Value = 4
$wsCages.Range("D43").Value = 5
$wsCages.Range("E43").Value = "wood"
$wsCages.Range("F43").Value = 8034
